{"js": "// Apply the prescription-form field updates described by the diff.\n// Each field value lives as the tail of a single run of text, e.g.\n// \"Patients Name: Swapnanil Bala\" -> \"Patients Name: Mahatab Ali\".\n// We locate each old full string via body.search() (exact, case-sensitive)\n// and replace it in place with the new full string.\n\nconst replacements = [\n  [\"Date&Time: 08-01-2025 21:00:59\", \"Date&Time: 16-01-2025 23:05:10\"],\n  [\"Patients Name: Swapnanil Bala\", \"Patients Name: Mahatab Ali\"],\n  [\"Age: 27\", \"Age: 10\"],\n  [\"Address: Boston, Massachuttes\", \"Address: Bharsala\"],\n  [\"Short History with Complaints: No\", \"Short History with Complaints: pain right groin for last 3 days with fever.\"],\n  [\"Important Clinical Findings: not\", \"Important Clinical Findings: Anterior hip point(right) tender\"],\n  [\"Investigation Advised: working\", \"Investigation Advised: TC DC ESR, Hb, CRP\"],\n  [\"Advice to follow: yet\", \"Advice to follow: Tablet CETIL 250 mg 1 tab BDPC X 5 days.\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the prescription-form field updates described by the diff.\n# Each field value lives as the tail of a single run of text, e.g.\n# \"Patients Name: Swapnanil Bala\" -> \"Patients Name: Mahatab Ali\".\n# Use Find/Replace (wdReplaceAll = 2) scoped to the whole document body\n# so each old full string is swapped for its new full string in place.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Date&Time: 08-01-2025 21:00:59\", \"Date&Time: 16-01-2025 23:05:10\"),\n    @(\"Patients Name: Swapnanil Bala\", \"Patients Name: Mahatab Ali\"),\n    @(\"Age: 27\", \"Age: 10\"),\n    @(\"Address: Boston, Massachuttes\", \"Address: Bharsala\"),\n    @(\"Short History with Complaints: No\", \"Short History with Complaints: pain right groin for last 3 days with fever.\"),\n    @(\"Important Clinical Findings: not\", \"Important Clinical Findings: Anterior hip point(right) tender\"),\n    @(\"Investigation Advised: working\", \"Investigation Advised: TC DC ESR, Hb, CRP\"),\n    @(\"Advice to follow: yet\", \"Advice to follow: Tablet CETIL 250 mg 1 tab BDPC X 5 days.\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
